$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bewertung")

# Edit 1: C3 - replace formula "=3*4.8" with "=17.64"
$ws.Range("C3").Formula = "=17.64"

# Edit 2: E3 - remove formula, set plain value 0.13
$ws.Range("E3").Value = 0.13

# Edit 3: E8 - clear cell style (reset to default/General) and set value to 8.82
$ws.Range("E8").Value = 8.82
$ws.Range("E8").ClearFormats()

# Edit 4: F23 - change value from 12 to 6
$ws.Range("F23").Value = 6

# Update the saved view state (selection + zoom) to match the end-of-session state
$ws.Range("J17").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 130
